$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 7.895
$ws.Range("B6").Value = 5.906000000000001
$ws.Range("B7").Value = 6.032000000000001
$ws.Range("D7").Value = -7.647
$ws.Range("D12").Value = -7.302
$ws.Range("D15").Value = -8.270999999999999
$ws.Range("B16").Value = 5.352
$ws.Range("B20").Value = 8.466000000000001
$ws.Range("D20").Value = -7.957000000000001
$ws.Range("D21").Value = -8.18
$ws.Range("D22").Value = -7.865
$ws.Range("D23").Value = -7.869999999999999
$ws.Range("B28").Value = 6.171
$ws.Range("B29").Value = 5.215999999999999
$ws.Range("D29").Value = -7.010000000000001
$ws.Range("B32").Value = 6.572
$ws.Range("D34").Value = -7.918000000000001
$ws.Range("B40").Value = 9.370000000000001
$ws.Range("D42").Value = -8.106999999999999
$ws.Range("D43").Value = -7.773000000000001
$ws.Range("D44").Value = -7.860000000000001
$ws.Range("D45").Value = -7.531000000000001
$ws.Range("B46").Value = 5.971
$ws.Range("D46").Value = -8.342000000000002
$ws.Range("D50").Value = -8.096
$ws.Range("B51").Value = 4.862
$ws.Range("D51").Value = -8.103
$ws.Range("B52").Value = 5.825
$ws.Range("B57").Value = 5.188999999999999
$ws.Range("B59").Value = 4.515
$ws.Range("B62").Value = 5.261
$ws.Range("B66").Value = 5.773
$ws.Range("D66").Value = -7.465999999999999
$ws.Range("D67").Value = -7.204000000000001
$ws.Range("B73").Value = 6.605000000000001
$ws.Range("B74").Value = 9.204000000000001
$ws.Range("D79").Value = -7.6
$ws.Range("D84").Value = -8.300000000000001
$ws.Range("B92").Value = 5.882
$ws.Range("D92").Value = -6.572999999999999
$ws.Range("D97").Value = -8.276
$ws.Range("B100").Value = 5.848999999999999
